$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.10057849630809
$ws.Range("C2").Value = 11.37698729203081
$ws.Range("E2").Value = 16.59688031244306
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("G2").Value = 3.587836969224761
$ws.Range("O2").Value = 16.1585145288316

$ws.Range("B3").Value = 14.26868350561348
$ws.Range("C3").Value = 10.76881909942569
$ws.Range("E3").Value = 15.64775420677554
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("G3").Value = 3.590295359403196
$ws.Range("O3").Value = 16.31841538226039

$ws.Range("B4").Value = 13.73221647078803
$ws.Range("C4").Value = 10.37592823587151
$ws.Range("E4").Value = 15.03971881126983
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("G4").Value = 3.591880113520168
$ws.Range("O4").Value = 16.42480075204937

$ws.Range("B5").Value = 13.50731571966918
$ws.Range("C5").Value = 10.21102872314305
$ws.Range("E5").Value = 14.78585217213869
$ws.Range("F5").Value = 15.008197319934
$ws.Range("G5").Value = 3.59254491630988
$ws.Range("O5").Value = 16.4701969270732

$ws.Range("B6").Value = 13.46959717079582
$ws.Range("C6").Value = 10.18336123015006
$ws.Range("E6").Value = 14.74333899549794
$ws.Range("F6").Value = 14.96433081551589
$ws.Range("G6").Value = 3.592656456062468
$ws.Range("O6").Value = 16.47785775762977

$ws.Range("B7").Value = 13.72920858024521
$ws.Range("C7").Value = 10.37372359881691
$ws.Range("E7").Value = 15.03631933268789
$ws.Range("F7").Value = 15.26647399323133
$ws.Range("G7").Value = 3.591889002250952
$ws.Range("O7").Value = 16.42540473394759

$ws.Range("B8").Value = 14.81916334307873
$ws.Range("C8").Value = 11.17138817656628
$ws.Range("E8").Value = 16.27501025964296
$ws.Range("F8").Value = 16.5399640634477
$ws.Range("G8").Value = 3.58866903726285
$ws.Range("O8").Value = 16.21193305027942

$ws.Range("B9").Value = 16.74688857192679
$ws.Range("C9").Value = 12.57754132270443
$ws.Range("E9").Value = 18.61360637468819
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 3.582948912490946
$ws.Range("O9").Value = 15.85934373915154

$ws.Range("B10").Value = 18.02929488752533
$ws.Range("C10").Value = 13.51065448849431
$ws.Range("E10").Value = 20.2726560391692
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 3.579104088670131
$ws.Range("O10").Value = 15.64184137010998

$ws.Range("B11").Value = 18.58276029904653
$ws.Range("C11").Value = 13.9129329313141
$ws.Range("E11").Value = 20.98503927773836
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 3.577431697676962
$ws.Range("O11").Value = 15.55220301588404

$ws.Range("B12").Value = 18.78799118730951
$ws.Range("C12").Value = 14.06204310601366
$ws.Range("E12").Value = 21.24875178547164
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 3.576809353632264
$ws.Range("O12").Value = 15.51962187782693

$ws.Range("B13").Value = 18.74398542996855
$ws.Range("C13").Value = 14.03007336879097
$ws.Range("E13").Value = 21.19222516878587
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 3.576942900424838
$ws.Range("O13").Value = 15.52657778795521

$ws.Range("B14").Value = 18.59973230001565
$ws.Range("C14").Value = 13.92526508118121
$ws.Range("E14").Value = 21.00685614722373
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 3.577380277900693
$ws.Range("O14").Value = 15.54949508551007

$ws.Range("B15").Value = 18.51080463835326
$ws.Range("C15").Value = 13.86064632152512
$ws.Range("E15").Value = 20.89252521126337
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 3.577649608816804
$ws.Range("O15").Value = 15.56371083924995

$ws.Range("B16").Value = 17.99251790636597
$ws.Range("C16").Value = 13.48391504258623
$ws.Range("E16").Value = 20.22525139082524
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 3.57921491983779
$ws.Range("O16").Value = 15.64788880291369

$ws.Range("B17").Value = 17.66686320518416
$ws.Range("C17").Value = 13.24709162589515
$ws.Range("E17").Value = 19.8050845660391
$ws.Range("F17").Value = 20.20408069617459
$ws.Range("G17").Value = 3.580194769119316
$ws.Range("O17").Value = 15.70192999430557

$ws.Range("B18").Value = 17.47674247176683
$ws.Range("C18").Value = 13.10878886926816
$ws.Range("E18").Value = 19.55942964652814
$ws.Range("F18").Value = 19.95656407809808
$ws.Range("G18").Value = 3.580765570061527
$ws.Range("O18").Value = 15.73388732272799

$ws.Range("B19").Value = 17.41188961919109
$ws.Range("C19").Value = 13.0616043833986
$ws.Range("E19").Value = 19.4755694308064
$ws.Range("F19").Value = 19.87204792380562
$ws.Range("G19").Value = 3.580960075099934
$ws.Range("O19").Value = 15.74485692765616

$ws.Range("B20").Value = 17.70182124502114
$ws.Range("C20").Value = 13.27251825765834
$ws.Range("E20").Value = 19.85022433053811
$ws.Range("F20").Value = 20.24955283636157
$ws.Range("G20").Value = 3.580089715986387
$ws.Range("O20").Value = 15.69608655510414

$ws.Range("B21").Value = 18.64222148475511
$ws.Range("C21").Value = 13.95613752610697
$ws.Range("E21").Value = 21.0614674822927
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 3.577251512770157
$ws.Range("O21").Value = 15.54272652787774

$ws.Range("B22").Value = 19.23142995009985
$ws.Range("C22").Value = 14.38412123717443
$ws.Range("E22").Value = 21.81782783699374
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 3.575460401269853
$ws.Range("O22").Value = 15.45045318695732

$ws.Range("B23").Value = 18.91929673146743
$ws.Range("C23").Value = 14.15742705757384
$ws.Range("E23").Value = 21.4173586030368
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 3.576410533693135
$ws.Range("O23").Value = 15.49896482909649

$ws.Range("B24").Value = 17.68602573610572
$ws.Range("C24").Value = 13.26102957123496
$ws.Range("E24").Value = 19.82982940237258
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 3.580137187225917
$ws.Range("O24").Value = 15.69872560875873

$ws.Range("B25").Value = 16.2485494308624
$ws.Range("C25").Value = 12.21448973395496
$ws.Range("E25").Value = 17.96487938914879
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 3.584433208672118
$ws.Range("O25").Value = 15.9475283000748
